# Add Moto3 and Moto2 qualifying and races for Great Britain (Silverstone)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MotoGP")

# Existing Great Britain block (before edit) occupies rows 30-32:
#   30: MotoGP Qualifying
#   31: MotoGP Sprint Race
#   32: MotoGP Race
#
# Target layout (rows 30-36):
#   30: MotoGP Qualifying   (unchanged)
#   31: Moto3  Qualifying   (NEW)
#   32: Moto2  Qualifying   (NEW)
#   33: MotoGP Sprint Race  (was row 31)
#   34: Moto3  Race         (NEW)
#   35: MotoGP Race         (was row 32)
#   36: Moto2  Race         (NEW)
#
# Insert 4 blank rows at the right spots (applied top-down so each Insert
# call uses the final target row index).
$ws.Rows("31:31").Insert()
$ws.Rows("32:32").Insert()
$ws.Rows("34:34").Insert()
$ws.Rows("36:36").Insert()

# --- Row 31: Moto3 Qualifying ---
$ws.Range("A31").Value2 = "Great Britain"
$ws.Range("B31").Value2 = "Moto3"
$ws.Range("C31").Value2 = "Qualifying"
$ws.Range("D31").Formula = "=VLOOKUP(A31,locations,4)"
$ws.Range("E31").Value2 = 45143
$ws.Range("F31").Value2 = 0.53472222222222221
$ws.Range("G31").Formula = "=E31"
$ws.Range("H31").Value2 = 0.5625
$ws.Range("I31").Formula = "=VLOOKUP(A31,locations,2)"
$ws.Range("J31").Formula = "=VLOOKUP(A31,locations,3)"

# --- Row 32: Moto2 Qualifying ---
$ws.Range("A32").Value2 = "Great Britain"
$ws.Range("B32").Value2 = "Moto2"
$ws.Range("C32").Value2 = "Qualifying"
$ws.Range("D32").Formula = "=VLOOKUP(A32,locations,4)"
$ws.Range("E32").Value2 = 45143
$ws.Range("F32").Value2 = 0.57291666666666663
$ws.Range("G32").Formula = "=E32"
$ws.Range("H32").Value2 = 0.60069444444444442
$ws.Range("I32").Formula = "=VLOOKUP(A32,locations,2)"
$ws.Range("J32").Formula = "=VLOOKUP(A32,locations,3)"

# --- Row 33: MotoGP Sprint Race (carried down from the old row 31, values unchanged) ---
$ws.Range("A33").Value2 = "Great Britain"
$ws.Range("B33").Value2 = "MotoGP"
$ws.Range("C33").Value2 = "Sprint Race"
$ws.Range("D33").Formula = "=VLOOKUP(A33,locations,4)"
$ws.Range("E33").Value2 = 45143
$ws.Range("F33").Value2 = 0.625
$ws.Range("G33").Formula = "=E33"
$ws.Range("H33").Value2 = 0.64583333333333337
$ws.Range("I33").Formula = "=VLOOKUP(A33,locations,2)"
$ws.Range("J33").Formula = "=VLOOKUP(A33,locations,3)"

# --- Row 34: Moto3 Race ---
$ws.Range("A34").Value2 = "Great Britain"
$ws.Range("B34").Value2 = "Moto3"
$ws.Range("C34").Value2 = "Race"
$ws.Range("D34").Formula = "=VLOOKUP(A34,locations,4)"
$ws.Range("E34").Value2 = 45144
$ws.Range("F34").Value2 = 0.46875
$ws.Range("G34").Formula = "=E34"
$ws.Range("H34").Value2 = 0.49305555555555558
$ws.Range("I34").Formula = "=VLOOKUP(A34,locations,2)"
$ws.Range("J34").Formula = "=VLOOKUP(A34,locations,3)"

# --- Row 35: MotoGP Race (carried down from the old row 32, values unchanged) ---
$ws.Range("A35").Value2 = "Great Britain"
$ws.Range("B35").Value2 = "MotoGP"
$ws.Range("C35").Value2 = "Race"
$ws.Range("D35").Formula = "=VLOOKUP(A35,locations,4)"
$ws.Range("E35").Value2 = 45144
$ws.Range("F35").Value2 = 0.54166666666666663
$ws.Range("G35").Formula = "=E35"
$ws.Range("H35").Value2 = 0.57638888888888895
$ws.Range("I35").Formula = "=VLOOKUP(A35,locations,2)"
$ws.Range("J35").Formula = "=VLOOKUP(A35,locations,3)"

# --- Row 36: Moto2 Race ---
$ws.Range("A36").Value2 = "Great Britain"
$ws.Range("B36").Value2 = "Moto2"
$ws.Range("C36").Value2 = "Race"
$ws.Range("D36").Formula = "=VLOOKUP(A36,locations,4)"
$ws.Range("E36").Value2 = 45144
$ws.Range("F36").Value2 = 0.60416666666666663
$ws.Range("G36").Formula = "=E36"
$ws.Range("H36").Value2 = 0.63194444444444442
$ws.Range("I36").Formula = "=VLOOKUP(A36,locations,2)"
$ws.Range("J36").Formula = "=VLOOKUP(A36,locations,3)"

# Restore the single-cell selection (Excel collapses multi-cell selections
# left over from editing back to the active cell on save).
$ws.Range("A1").Select()
